$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed / new cell values (April -> May 2023 dark-sky times) ---
$ws.Range("B3").Value = "2023/05/01 04:46"
$ws.Range("C3").Value = "2023/05/01 22:10"
$ws.Range("D3").Value = "2023/05/01 04:41"
$ws.Range("E3").Value = "2023/05/01 16:25"
$ws.Range("B4").Value = "2023/05/02 05:06"
$ws.Range("C4").Value = "2023/05/02 22:11"
$ws.Range("D4").Value = "2023/05/02 04:40"
$ws.Range("E4").Value = "2023/05/02 17:28"
$ws.Range("B5").Value = "2023/05/03 05:27"
$ws.Range("C5").Value = "2023/05/03 22:13"
$ws.Range("D5").Value = "2023/05/03 04:38"
$ws.Range("E5").Value = "2023/05/03 18:33"
$ws.Range("B6").Value = "2023/05/04 05:49"
$ws.Range("C6").Value = "2023/05/04 22:14"
$ws.Range("D6").Value = "2023/05/04 04:36"
$ws.Range("E6").Value = "2023/05/04 19:41"
$ws.Range("B7").Value = "2023/05/05 06:15"
$ws.Range("C7").Value = "2023/05/05 22:16"
$ws.Range("D7").Value = "2023/05/05 04:34"
$ws.Range("E7").Value = "2023/05/05 20:53"
$ws.Range("B8").Value = "2023/05/06 06:45"
$ws.Range("C8").Value = "2023/05/06 22:18"
$ws.Range("D8").Value = "2023/05/06 04:33"
$ws.Range("E8").Value = "2023/05/06 22:07"
$ws.Range("B9").Value = "2023/05/07 07:23"
$ws.Range("C9").Value = "2023/05/07 22:19"
$ws.Range("D9").Value = "2023/05/07 04:31"
$ws.Range("E9").Value = "2023/05/07 23:20"
$ws.Range("B10").Value = "2023/05/08 08:12"
$ws.Range("C10").Value = "2023/05/08 22:21"
$ws.Range("D10").Value = "2023/05/08 04:29"
$ws.Range("B11").Value = "2023/05/09 09:12"
$ws.Range("C11").Value = "2023/05/09 22:22"
$ws.Range("D11").Value = "2023/05/09 04:28"
$ws.Range("E11").Value = "2023/05/09 00:29"
$ws.Range("E11").Font.Bold = $false
$ws.Range("E11").HorizontalAlignment = 7
$ws.Range("E11").WrapText = $true
$ws.Range("B12").Value = "2023/05/10 10:22"
$ws.Range("C12").Value = "2023/05/10 22:24"
$ws.Range("D12").Value = "2023/05/10 04:26"
$ws.Range("E12").Value = "2023/05/10 01:27"
$ws.Range("B13").Value = "2023/05/11 11:37"
$ws.Range("C13").Value = "2023/05/11 22:25"
$ws.Range("D13").Value = "2023/05/11 04:24"
$ws.Range("E13").Value = "2023/05/11 02:15"
$ws.Range("B14").Value = "2023/05/12 12:54"
$ws.Range("C14").Value = "2023/05/12 22:27"
$ws.Range("D14").Value = "2023/05/12 04:23"
$ws.Range("E14").Value = "2023/05/12 02:53"
$ws.Range("B15").Value = "2023/05/13 14:09"
$ws.Range("C15").Value = "2023/05/13 22:28"
$ws.Range("D15").Value = "2023/05/13 04:21"
$ws.Range("E15").Value = "2023/05/13 03:23"
$ws.Range("B16").Value = "2023/05/14 15:22"
$ws.Range("C16").Value = "2023/05/14 22:30"
$ws.Range("D16").Value = "2023/05/14 04:20"
$ws.Range("E16").Value = "2023/05/14 03:50"
$ws.Range("B17").Value = "2023/05/15 16:33"
$ws.Range("C17").Value = "2023/05/15 22:32"
$ws.Range("D17").Value = "2023/05/15 04:18"
$ws.Range("E17").Value = "2023/05/15 04:13"
$ws.Range("B18").Value = "2023/05/16 17:43"
$ws.Range("C18").Value = "2023/05/16 22:33"
$ws.Range("D18").Value = "2023/05/16 04:17"
$ws.Range("E18").Value = "2023/05/16 04:36"
$ws.Range("B19").Value = "2023/05/17 18:54"
$ws.Range("C19").Value = "2023/05/17 22:35"
$ws.Range("D19").Value = "2023/05/17 04:15"
$ws.Range("E19").Value = "2023/05/17 05:00"
$ws.Range("B20").Value = "2023/05/18 20:05"
$ws.Range("C20").Value = "2023/05/18 22:36"
$ws.Range("D20").Value = "2023/05/18 04:14"
$ws.Range("E20").Value = "2023/05/18 05:26"
$ws.Range("B21").Value = "2023/05/19 21:14"
$ws.Range("C21").Value = "2023/05/19 22:38"
$ws.Range("D21").Value = "2023/05/19 04:12"
$ws.Range("E21").Value = "2023/05/19 05:56"
$ws.Range("B22").Value = "2023/05/20 22:21"
$ws.Range("C22").Value = "2023/05/20 22:39"
$ws.Range("D22").Value = "2023/05/20 04:11"
$ws.Range("E22").Value = "2023/05/20 06:32"
$ws.Range("B23").Value = "2023/05/21 23:22"
$ws.Range("C23").Value = "2023/05/21 22:41"
$ws.Range("D23").Value = "2023/05/21 04:10"
$ws.Range("E23").Value = "2023/05/21 07:14"
$ws.Range("C24").Value = "2023/05/22 22:42"
$ws.Range("D24").Value = "2023/05/22 04:08"
$ws.Range("E24").Value = "2023/05/22 08:04"
$ws.Range("B25").Value = "2023/05/23 00:14"
$ws.Range("B25").Font.Bold = $false
$ws.Range("B25").HorizontalAlignment = 7
$ws.Range("B25").WrapText = $true
$ws.Range("C25").Value = "2023/05/23 22:43"
$ws.Range("D25").Value = "2023/05/23 04:07"
$ws.Range("E25").Value = "2023/05/23 09:01"
$ws.Range("B26").Value = "2023/05/24 00:58"
$ws.Range("C26").Value = "2023/05/24 22:45"
$ws.Range("D26").Value = "2023/05/24 04:06"
$ws.Range("E26").Value = "2023/05/24 10:02"
$ws.Range("B27").Value = "2023/05/25 01:33"
$ws.Range("C27").Value = "2023/05/25 22:46"
$ws.Range("D27").Value = "2023/05/25 04:05"
$ws.Range("E27").Value = "2023/05/25 11:04"
$ws.Range("B28").Value = "2023/05/26 02:02"
$ws.Range("C28").Value = "2023/05/26 22:48"
$ws.Range("D28").Value = "2023/05/26 04:03"
$ws.Range("E28").Value = "2023/05/26 12:07"
$ws.Range("B29").Value = "2023/05/27 02:27"
$ws.Range("C29").Value = "2023/05/27 22:49"
$ws.Range("D29").Value = "2023/05/27 04:02"
$ws.Range("E29").Value = "2023/05/27 13:09"
$ws.Range("B30").Value = "2023/05/28 02:49"
$ws.Range("C30").Value = "2023/05/28 22:50"
$ws.Range("D30").Value = "2023/05/28 04:01"
$ws.Range("E30").Value = "2023/05/28 14:10"
$ws.Range("B31").Value = "2023/05/29 03:09"
$ws.Range("C31").Value = "2023/05/29 22:52"
$ws.Range("D31").Value = "2023/05/29 04:00"
$ws.Range("E31").Value = "2023/05/29 15:12"
$ws.Range("B32").Value = "2023/05/30 03:29"
$ws.Range("C32").Value = "2023/05/30 22:53"
$ws.Range("D32").Value = "2023/05/30 03:59"
$ws.Range("E32").Value = "2023/05/31 17:22"
$ws.Range("E32").Font.Bold = $false
$ws.Range("E32").HorizontalAlignment = 7
$ws.Range("E32").WrapText = $true
$ws.Range("B33").Value = "2023/06/01 04:14"
$ws.Range("B33").Font.Bold = $false
$ws.Range("B33").HorizontalAlignment = 7
$ws.Range("B33").WrapText = $true
$ws.Range("C33").Value = "2023/05/31 22:54"
$ws.Range("C33").Font.Bold = $false
$ws.Range("C33").HorizontalAlignment = 7
$ws.Range("C33").WrapText = $true
$ws.Range("D33").Value = "2023/05/31 03:58"
$ws.Range("D33").Font.Bold = $false
$ws.Range("D33").HorizontalAlignment = 7
$ws.Range("D33").WrapText = $true
$ws.Range("A34").Value = 32
$ws.Range("A34").Font.Bold = $true
$ws.Range("A34").HorizontalAlignment = 7
$ws.Range("A34").WrapText = $true

# --- Clear cells that no longer have data ---
$ws.Range("E10").ClearContents() | Out-Null
$ws.Range("B24").ClearContents() | Out-Null
